# resource_log.xlsx — refresh of the Power-Query-backed "resource_log" table.
# The underlying resource_log.csv now only contains a single "process" run
# (the earlier "clear" and "train" stages are gone), so the query table
# shrinks from 3 data rows to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("resource_log")

# Drop the two stale rows ("clear" in row 2, "train" in row 3). Deleting the
# entire sheet row (not just the ListRow) shifts the remaining data up and
# shrinks the table/used-range automatically, just like a real query refresh
# that returns fewer records.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# The sole remaining data row (now row 2) becomes the fresh "process" entry
# returned by the refreshed query.
$ws.Range("A2").Value = "process"
$ws.Range("B2").Value = 45790.785466111112
$ws.Range("C2").Value = 45790.787230532405
$ws.Range("D2").Value = 152.445155
$ws.Range("E2").Value = 3.8
$ws.Range("F2").Value = 3.4
$ws.Range("G2").Value = 52.8
$ws.Range("H2").Value = 52.7

# A query refresh writes plain (General-formatted) cells for the new row,
# rather than carrying over the old explicit "General" style - so drop the
# leftover formatting on the non-date columns.
$ws.Range("A2").ClearFormats()
$ws.Range("D2:L2").ClearFormats()

# gpu_* columns (I:L) come back blank for this run - clear them out entirely
# rather than leaving empty-but-present cells.
$ws.Range("I2:L2").ClearContents()

# The hidden ExternalData_1 name tracks the query's spill range; shrink it
# to match the new 1-row result (header + 1 data row, columns A:L).
$dn = $wb.Names.Item("ExternalData_1")
$dn.RefersTo = "=resource_log!`$A`$1:`$L`$2"

# Selection follows the shrunk calculated column (was M2:M4, now just M2).
$ws.Range("M2").Select()
